$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Order Details")

$ws.Range("C2").Value = 717
$ws.Range("C3").Value = 15
$ws.Range("C4").Value = 494
$ws.Range("C5").Value = 591
$ws.Range("C6").Value = 698
$ws.Range("C7").Value = 650
$ws.Range("C8").Value = 293
$ws.Range("C9").Value = 268
$ws.Range("C10").Value = 476
$ws.Range("C11").Value = 410
$ws.Range("C12").Value = 499
$ws.Range("C13").Value = 395
$ws.Range("C14").Value = 67
$ws.Range("C15").Value = 871
$ws.Range("C16").Value = 898
$ws.Range("C18").Value = 14
$ws.Range("C19").Value = 777
$ws.Range("C20").Value = 712
$ws.Range("C21").Value = 159
$ws.Range("C23").Value = 146
$ws.Range("C24").Value = 571
$ws.Range("C25").Value = 717
$ws.Range("C26").Value = 295
$ws.Range("C27").Value = 989
$ws.Range("C28").Value = 227
$ws.Range("C29").Value = 549
$ws.Range("C30").Value = 104
$ws.Range("C31").Value = 729
$ws.Range("C32").Value = 96
$ws.Range("C33").Value = 928
$ws.Range("C34").Value = 376
$ws.Range("C35").Value = 401
$ws.Range("C36").Value = 911
$ws.Range("C37").Value = 781
$ws.Range("C38").Value = 774
$ws.Range("C39").Value = 905
$ws.Range("C40").Value = 231
$ws.Range("C41").Value = 32
$ws.Range("C42").Value = 191
$ws.Range("C43").Value = 191
$ws.Range("C44").Value = 744
$ws.Range("C45").Value = 298
$ws.Range("C46").Value = 672
$ws.Range("C47").Value = 74
$ws.Range("C48").Value = 431
$ws.Range("C49").Value = 106
$ws.Range("C50").Value = 710
$ws.Range("C51").Value = 69
$ws.Range("C52").Value = 299
$ws.Range("C53").Value = 700
$ws.Range("C54").Value = 804
$ws.Range("C55").Value = 143
$ws.Range("C56").Value = 388
$ws.Range("C57").Value = 372
$ws.Range("C58").Value = 211
$ws.Range("C59").Value = 825
$ws.Range("C60").Value = 276
$ws.Range("C61").Value = 793
$ws.Range("C62").Value = 390
$ws.Range("C63").Value = 527
$ws.Range("C64").Value = 249
$ws.Range("C65").Value = 435
$ws.Range("C66").Value = 370
$ws.Range("C67").Value = 346
$ws.Range("C68").Value = 46
$ws.Range("C70").Value = 397
